# Converter gibt XML-Temp-Datei nun aus
# Update the sample/template data of the ILIAS course-import sheet:
# - New "Kurs-Referenz" ids (column A) for the group rows
# - Refreshed hierarchy ids (column D)
# - Course admins (column F) unified to "root"
# - Obsolete "crs | grp" / "Typ" column (P) cleared out
# - Selection cursor moved to A11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row P1/P2: drop the leftover "crs | grp" / "Typ" labels ---
$ws.Range("P1").ClearContents()
$ws.Range("P2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = 66
$ws.Range("D3").Value = 65

# --- Row 4 ---
$ws.Range("A4").Value = 112
$ws.Range("D4").Value = 66
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = "root"
$ws.Range("P4").ClearContents()

# --- Row 5 ---
$ws.Range("A5").Value = 113
$ws.Range("D5").Value = 66
$ws.Range("F3").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = "root"
$ws.Range("P5").ClearContents()

# --- Row 6 ---
$ws.Range("A6").Value = 114
$ws.Range("D6").Value = 66
$ws.Range("F3").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Value = "root"
$ws.Range("P6").ClearContents()

# --- Row 7 ---
$ws.Range("A7").Value = 115
$ws.Range("D7").Value = 66
$ws.Range("F3").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = "root"
$ws.Range("P7").ClearContents()

# --- Row 8 ---
$ws.Range("A8").Value = 116
$ws.Range("D8").Value = 66
$ws.Range("F3").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = "root"
$ws.Range("P8").ClearContents()

$excel.CutCopyMode = $false

# --- Move the selection cursor ---
$ws.Range("A11").Select()
